$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.384.65"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "2.621.10"
$ws.Range("E3").Value = "  +8.48%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'313.85"
$ws.Range("E5").Value = "  +4.38%  "
$ws.Range("D6").Value = "'101.50"
$ws.Range("E6").Value = "  +4.55%  "
$ws.Range("D7").Value = "'0.602"
$ws.Range("E7").Value = "  +6.35%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = "  +14.38%  "
$ws.Range("D10").Value = "'38.91"
$ws.Range("E10").Value = "  +10.75%  "
$ws.Range("D11").Value = "'54.56"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "'0.0843"
$ws.Range("E12").Value = "  +6.17%  "
$ws.Range("D13").Value = "'8.36"
$ws.Range("E13").Value = "  +16.49%  "
$ws.Range("D14").Value = "3.012.55"
$ws.Range("E14").Value = "  +8.01%  "
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "2.615.41"
$ws.Range("E16").Value = "  +8.91%  "
$ws.Range("D17").Value = "'0.917"
$ws.Range("E17").Value = "  +7.99%  "
$ws.Range("D18").Value = "'15.19"
$ws.Range("E18").Value = "  +6.05%  "
$ws.Range("D19").Value = "46.555.64"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").Value = "'13.41"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "'0.0000103"
$ws.Range("E21").Value = "  +8.01%  "
$ws.Range("D22").Value = "'6.79"
$ws.Range("E22").Value = "  +8.68%  "
$ws.Range("D23").Value = "'71.06"
$ws.Range("E23").Value = "  +5.33%  "
$ws.Range("D24").Value = "'255.66"
$ws.Range("E24").Value = "  +4.93%  "
$ws.Range("D25").Value = "'3.10"
$ws.Range("E25").Value = "  +9.93%  "
$ws.Range("D26").Value = "'2.22"
$ws.Range("E26").Value = "  +13.96%  "
$ws.Range("D27").Value = "'28.15"
$ws.Range("E27").Value = "  +31.11%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'10.71"
$ws.Range("E29").Value = "  +8.98%  "
$ws.Range("D30").Value = "'41.13"
$ws.Range("E30").Value = "  +5.62%  "
$ws.Range("D31").Value = "'2.29"
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("D32").Value = "'6.23"
$ws.Range("E32").Value = "  +11.53%  "
$ws.Range("D33").Value = "'3.75"
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("D34").Value = "'2.29"
$ws.Range("E34").Value = "  +14.05%  "
$ws.Range("D35").Value = "'2.86"
$ws.Range("E35").Value = "  +4.13%  "
$ws.Range("D36").Value = "'154.24"
$ws.Range("E36").Value = "  +3.90%  "
$ws.Range("D37").Value = "'0.0841"
$ws.Range("E37").Value = "  +7.93%  "
$ws.Range("E38").Value = "  +5.00%  "
$ws.Range("E39").Value = "  +6.29%  "
$ws.Range("D40").Value = "'17.08"
$ws.Range("E40").Value = "  +11.46%  "
$ws.Range("D41").Value = "'4.27"
$ws.Range("E41").Value = "  +8.99%  "
$ws.Range("D42").Value = "'3.65"
$ws.Range("E42").Value = "  +11.15%  "
$ws.Range("D43").Value = "'0.0329"
$ws.Range("E43").Value = "  +9.08%  "
$ws.Range("D44").Value = "'21.04"
$ws.Range("E44").Value = "  +36.57%  "
$ws.Range("D45").Value = "2.034.33"
$ws.Range("E45").Value = "  +4.18%  "
$ws.Range("D46").Value = "'0.998"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "'91.62"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").Value = "'113.12"
$ws.Range("E48").Value = "  +9.90%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.82"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'9.24"
$ws.Range("E50").Value = "  +6.65%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'78.39"
$ws.Range("E51").Value = "  +14.07%  "
